$d = $word.ActiveDocument

# --- Edit 1: swap the names of the two coincident ToC bookmarks that ---
# --- wrap the "Заявление о предоставлении..." title run.             ---
$bmA = $d.Bookmarks.Item("_Toc405368347")
$bmStart = $bmA.Start
$bmEnd = $bmA.End
$bmB = $d.Bookmarks.Item("_Toc373232731")
$bmA.Delete()
$bmB.Delete()
$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_Toc373232731", $bmRange)
$d.Bookmarks.Add("_Toc405368347", $bmRange)

# --- Edit 2: the "spouse is not an employee" paragraph becomes bold   ---
# --- and its literal sentence is replaced with the ${IS_DZO} merge   ---
# --- placeholder.                                                    ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Супруг не является работником Общества или ДЗО.`r") {
        $target = $p
        break
    }
}

$r = $target.Range
$r.Font.Bold = $true
$r.Font.BoldBi = $true
$r.Text = '${IS_DZO}'
